# Data_Mapping.xlsx - "Add files via upload" edit
# Updates the "Data Mapping" worksheet:
#   - Destination Column for CAMIS  -> resturant_id
#   - Destination Column for DBA    -> resturant_name
#   - Destination Column for BORO   -> borough
#   - Data Type for PHONE           -> Text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Mapping")

$ws.Range("E3").Value = "resturant_id"
$ws.Range("E4").Value = "resturant_name"
$ws.Range("E5").Value = "borough"
$ws.Range("C9").Value = "Text"

$ws.Range("D5:D12").Select()
